$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 94 (pushes the existing "note" row 94 down to 95),
# inheriting formatting from the row above (row 93).
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row with the next day's data.
$ws.Range("A94").Value = 43949
$ws.Range("B94").Value = 396
$ws.Range("C94").Value = 31114
$ws.Range("D94").Value = 148
$ws.Range("E94").Value = 6664

# Update the selected/active cell shown in the saved view.
$ws.Range("E96").Select()

# Extend the print area to include the newly added row.
$n = $wb.Names.Item(1)
$n.RefersTo = "=相談件数!`$A`$1:`$E`$97"
